$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 20 new traceability rows (id 30-49 / sheet rows 31-50) to the "data" sheet.
# initial_date / final_date columns (E, F) reuse the existing "YYYY-MM-DD HH:MM:SS" date-time
# style already used throughout the sheet, except for the last row which introduces a new
# date-only "YYYY-MM-DD" style.

# Row 31
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "check.flac"
$ws.Range("C31").Value = "Heather Carrillo (@31725671)"
$ws.Range("D31").Value = "Robin Blair (@16670044)"
$ws.Range("E31").Value = 45374
$ws.Range("E31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F31").Value = 45377
$ws.Range("F31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = "Skill pick find cell trouble method ten. Food later baby both. Ago us girl find. Laugh discover contain."
$ws.Range("I31").Value = "Jamie Willis"
$ws.Range("J31").Value = "RECHAZADO - CENCO"
$ws.Range("K31").Value = "Type 1"
$ws.Range("L31").Value = "Jon Palmer (@23400841),Heather Carrillo (@31725671),Jon Palmer (@23400841)"

# Row 32
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "sound.flac"
$ws.Range("C32").Value = "Jennifer Fuentes (@9840607)"
$ws.Range("D32").Value = "Jon Palmer (@23400841)"
$ws.Range("E32").Value = 45380
$ws.Range("E32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F32").Value = 45399
$ws.Range("F32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G32").Value = -1
$ws.Range("H32").Value = "Exactly activity agreement store program seem politics across. Ten commercial employee senior democratic region local. Benefit why collection feeling."
$ws.Range("I32").Value = "Nathan Brown"
$ws.Range("J32").Value = "RECHAZADO - CENCO"
$ws.Range("K32").Value = "Type 2"
$ws.Range("L32").Value = "Jon Palmer (@23400841),Jon Palmer (@23400841),Heather Carrillo (@31725671)"

# Row 33
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "half.mov"
$ws.Range("C33").Value = "Heather Carrillo (@31725671)"
$ws.Range("D33").Value = "Melissa Cox (@51800599)"
$ws.Range("E33").Value = 45368
$ws.Range("E33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F33").Value = 45382
$ws.Range("F33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G33").Value = 11
$ws.Range("H33").Value = "Gas help region those walk everybody seek hair. Account measure such popular part. Book issue fill tough natural rest. Seven few necessary model direction contain.`nFirst about song begin."
$ws.Range("I33").Value = "Nicole Hunter"
$ws.Range("J33").Value = "RECHAZADO - DECANO"
$ws.Range("K33").Value = "Type 3"
$ws.Range("L33").Value = "Heather Carrillo (@31725671),Kim Rosales (@39851784)"

# Row 34
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "commercial.xlsx"
$ws.Range("C34").Value = "Kevin Foster (@32273703)"
$ws.Range("D34").Value = "Kevin Foster (@32273703)"
$ws.Range("E34").Value = 45355
$ws.Range("E34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F34").Value = 45356
$ws.Range("F34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G34").Value = 24
$ws.Range("H34").Value = "Hear current power field. Until little environmental clearly though decision. Win fast address million push financial. Huge best dinner situation."
$ws.Range("I34").Value = "Kimberly Hudson"
$ws.Range("J34").Value = "CERRADO"
$ws.Range("K34").Value = "Type 3"
$ws.Range("L34").Value = "Heather Carrillo (@31725671),Jon Palmer (@23400841),Kim Rosales (@39851784)"

# Row 35
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "page.jpeg"
$ws.Range("C35").Value = "Chad Fox (@44758173)"
$ws.Range("D35").Value = "Kim Rosales (@39851784)"
$ws.Range("E35").Value = 45359
$ws.Range("E35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F35").Value = 45386
$ws.Range("F35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G35").Value = 20
$ws.Range("H35").Value = "Walk tough industry pass radio world travel. Would finish PM Mr billion agency right. Camera exactly class identify."
$ws.Range("I35").Value = "Brooke Oconnell"
$ws.Range("J35").Value = "CERRADO"
$ws.Range("K35").Value = "Type 1"
$ws.Range("L35").Value = "Kim Rosales (@39851784),Jon Palmer (@23400841),Jon Palmer (@23400841)"

# Row 36
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "TV.xlsx"
$ws.Range("C36").Value = "Kevin Foster (@32273703)"
$ws.Range("D36").Value = "Robin Blair (@16670044)"
$ws.Range("E36").Value = 45369
$ws.Range("E36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F36").Value = 45379
$ws.Range("F36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G36").Value = 10
$ws.Range("H36").Value = "However quickly partner two yeah production into mother. Carry account near move this record. Third activity oil short.`nLaugh animal responsibility from technology."
$ws.Range("I36").Value = "Kristy Carney"
$ws.Range("J36").Value = "APROBADO - CENCO"
$ws.Range("K36").Value = "Type 2"
$ws.Range("L36").Value = "Jon Palmer (@23400841),Jon Palmer (@23400841)"

# Row 37
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "leg.xlsx"
$ws.Range("C37").Value = "Chad Fox (@44758173)"
$ws.Range("D37").Value = "Heather Carrillo (@31725671)"
$ws.Range("E37").Value = 45353
$ws.Range("E37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F37").Value = 45374
$ws.Range("F37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G37").Value = 26
$ws.Range("H37").Value = "Remember agreement later every sort south. Anything ready off research still night paper. Impact PM letter money since."
$ws.Range("I37").Value = "Jared Tucker"
$ws.Range("J37").Value = "RECHAZADO - DECANO"
$ws.Range("K37").Value = "Type 1"
$ws.Range("L37").Value = "Kim Rosales (@39851784)"

# Row 38
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "month.docx"
$ws.Range("C38").Value = "Melissa Cox (@51800599)"
$ws.Range("D38").Value = "Desiree Robinson (@58036467)"
$ws.Range("E38").Value = 45372
$ws.Range("E38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F38").Value = 45377
$ws.Range("F38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G38").Value = 7
$ws.Range("H38").Value = "Across new debate discussion least everybody. Decision reflect real painting throw operation."
$ws.Range("I38").Value = "Michael Huff"
$ws.Range("J38").Value = "APROBADO - DECANO"
$ws.Range("K38").Value = "Type 2"
$ws.Range("L38").Value = "Heather Carrillo (@31725671)"

# Row 39
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "week.jpeg"
$ws.Range("C39").Value = "Kevin Foster (@32273703)"
$ws.Range("D39").Value = "Heather Carrillo (@31725671)"
$ws.Range("E39").Value = 45382
$ws.Range("E39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F39").Value = 45383
$ws.Range("F39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G39").Value = -3
$ws.Range("H39").Value = "Anything yet hope. Plan sing traditional fall. Finally even hotel agent discuss.`nRealize field wait simple.`nMost argue add protect fill business give. Commercial force never past."
$ws.Range("I39").Value = "Willie Murphy"
$ws.Range("J39").Value = "EN PROCESO"
$ws.Range("K39").Value = "Type 3"
$ws.Range("L39").Value = "Jon Palmer (@23400841),Jon Palmer (@23400841)"

# Row 40
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "large.pages"
$ws.Range("C40").Value = "Desiree Robinson (@58036467)"
$ws.Range("D40").Value = "Chad Fox (@44758173)"
$ws.Range("E40").Value = 45349
$ws.Range("E40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F40").Value = 45357
$ws.Range("F40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G40").Value = 30
$ws.Range("H40").Value = "Collection total form family reduce power bank. Wide this culture production magazine drop.`nParty throw already its one expect until. Build drug certain always difference difference cut."
$ws.Range("I40").Value = "Sarah Roberts"
$ws.Range("J40").Value = "CERRADO"
$ws.Range("K40").Value = "Type 3"
$ws.Range("L40").Value = "Jon Palmer (@23400841),Jon Palmer (@23400841)"

# Row 41
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "next.css"
$ws.Range("C41").Value = "Richard Thomas (@57640453)"
$ws.Range("D41").Value = "Melissa Johnson (@47666474)"
$ws.Range("E41").Value = 45355
$ws.Range("E41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F41").Value = 45370
$ws.Range("F41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G41").Value = 24
$ws.Range("H41").Value = "Movie young chair she off. Father serious painting positive voice become generation.`nPage ground too make everybody. Week subject able this back you. West international exist stand why authority one."
$ws.Range("I41").Value = "Crystal Smith"
$ws.Range("J41").Value = "RECHAZADO - CONTABILIDAD"
$ws.Range("K41").Value = "Type 3"
$ws.Range("L41").Value = "Melissa Johnson (@47666474)"

# Row 42
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "apply.mp3"
$ws.Range("C42").Value = "Melissa Johnson (@47666474)"
$ws.Range("D42").Value = "Benjamin Miller (@1477097)"
$ws.Range("E42").Value = 45381
$ws.Range("E42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F42").Value = 45408
$ws.Range("F42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G42").Value = -2
$ws.Range("H42").Value = "Thought military task another billion. Carry audience trouble apply.`nEducation relate reach turn tell discover unit. Number only room ask.`nInterview record western scientist."
$ws.Range("I42").Value = "Jason Palmer"
$ws.Range("J42").Value = "RECHAZADO - CENCO"
$ws.Range("K42").Value = "Type 2"
$ws.Range("L42").Value = "Melissa Johnson (@47666474),Taylor Wilkins (@5464236),Taylor Wilkins (@5464236)"

# Row 43
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "author.js"
$ws.Range("C43").Value = "Melissa Johnson (@47666474)"
$ws.Range("D43").Value = "Melissa Johnson (@47666474)"
$ws.Range("E43").Value = 45381
$ws.Range("E43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F43").Value = 45382
$ws.Range("F43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G43").Value = -2
$ws.Range("H43").Value = "Outside important we pull reveal early. Finally house crime region again poor.`nMagazine move who. Last hair human idea. Without part through president."
$ws.Range("I43").Value = "Melissa Johnson"
$ws.Range("J43").Value = "RECHAZADO - CENCO"
$ws.Range("K43").Value = "Type 2"
$ws.Range("L43").Value = "Kenneth Chang (@59512309),Maria Lewis (@48114355)"

# Row 44
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "case.mov"
$ws.Range("C44").Value = "Benjamin Miller (@1477097)"
$ws.Range("D44").Value = "Richard Thomas (@57640453)"
$ws.Range("E44").Value = 45360
$ws.Range("E44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F44").Value = 45376
$ws.Range("F44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G44").Value = 19
$ws.Range("H44").Value = "Author pick including able surface. Amount crime population develop clear late.`nChoose form argue more indicate contain. Pick example often entire. Course drive song good method."
$ws.Range("I44").Value = "Rachel Barton"
$ws.Range("J44").Value = "APROBADO - CENCO"
$ws.Range("K44").Value = "Type 3"
$ws.Range("L44").Value = "Kenneth Chang (@59512309),Taylor Wilkins (@5464236),Taylor Wilkins (@5464236)"

# Row 45
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "red.gif"
$ws.Range("C45").Value = "Shannon Brown (@90416220)"
$ws.Range("D45").Value = "Maria Lewis (@48114355)"
$ws.Range("E45").Value = 45359
$ws.Range("E45").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F45").Value = 45360
$ws.Range("F45").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G45").Value = 20
$ws.Range("H45").Value = "Network second serve arm. Full thank son send.`nDiscover court after. Article yard when write music forget.`nMinute especially better people final. Third its nearly he. Hold scientist magazine walk."
$ws.Range("I45").Value = "Ryan Delacruz"
$ws.Range("J45").Value = "RECHAZADO - DECANO"
$ws.Range("K45").Value = "Type 1"
$ws.Range("L45").Value = "Melissa Johnson (@47666474)"

# Row 46
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "Mrs.html"
$ws.Range("C46").Value = "Richard Thomas (@57640453)"
$ws.Range("D46").Value = "Benjamin Miller (@1477097)"
$ws.Range("E46").Value = 45353
$ws.Range("E46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F46").Value = 45360
$ws.Range("F46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G46").Value = 26
$ws.Range("H46").Value = "Data mean truth despite describe step arrive seven. A again job television.`nYet smile forward like nice attack these. Mother card data wrong hand."
$ws.Range("I46").Value = "Vincent Guerra"
$ws.Range("J46").Value = "RECHAZADO - CENCO"
$ws.Range("K46").Value = "Type 1"
$ws.Range("L46").Value = "Taylor Wilkins (@5464236),Maria Lewis (@48114355)"

# Row 47
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "tend.mp4"
$ws.Range("C47").Value = "Melissa Johnson (@47666474)"
$ws.Range("D47").Value = "Maria Lewis (@48114355)"
$ws.Range("E47").Value = 45373
$ws.Range("E47").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F47").Value = 45383
$ws.Range("F47").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G47").Value = 6
$ws.Range("H47").Value = "Future political simple star. Seem see join vote condition continue. Lead phone bad collection."
$ws.Range("I47").Value = "Michael Dunn"
$ws.Range("J47").Value = "CERRADO"
$ws.Range("K47").Value = "Type 2"
$ws.Range("L47").Value = "Melissa Johnson (@47666474),Kenneth Chang (@59512309)"

# Row 48
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "last.pdf"
$ws.Range("C48").Value = "Julia Herman (@87635666)"
$ws.Range("D48").Value = "Julia Herman (@87635666)"
$ws.Range("E48").Value = 45377
$ws.Range("E48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F48").Value = 45406
$ws.Range("F48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = "Indicate your series describe. Such imagine few control most. Economic behind security especially whole each rise fast."
$ws.Range("I48").Value = "Katelyn Harrell"
$ws.Range("J48").Value = "EN PROCESO"
$ws.Range("K48").Value = "Type 2"
$ws.Range("L48").Value = "Maria Lewis (@48114355)"

# Row 49
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "structure.mp3"
$ws.Range("C49").Value = "Shannon Brown (@90416220)"
$ws.Range("D49").Value = "Maria Lewis (@48114355)"
$ws.Range("E49").Value = 45360
$ws.Range("E49").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F49").Value = 45376
$ws.Range("F49").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G49").Value = 19
$ws.Range("H49").Value = "Structure ask region voice professional really cup. Suggest too hour economy career right.`nPerform TV else color. Box nation cut wonder often speak effort beautiful. Expect high each get baby do."
$ws.Range("I49").Value = "Michael Johnston"
$ws.Range("J49").Value = "APROBADO - CENCO"
$ws.Range("K49").Value = "Type 1"
$ws.Range("L49").Value = "Kenneth Chang (@59512309),Taylor Wilkins (@5464236)"

# Row 50
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "program.flac"
$ws.Range("C50").Value = "Shannon Brown (@90416220)"
$ws.Range("D50").Value = "Kenneth Chang (@59512309)"
$ws.Range("E50").Value = 45360
# Register the lowercase numFmt variant first (left unused), then switch to the uppercase one actually applied -
# mirrors the existing unused 164 "yyyy-mm-dd h:mm:ss" / used 165 "YYYY-MM-DD HH:MM:SS" pairing already in the workbook.
$ws.Range("E50").NumberFormat = "yyyy-mm-dd"
$ws.Range("E50").NumberFormat = "YYYY-MM-DD"
$ws.Range("F50").Value = 45380
$ws.Range("F50").NumberFormat = "YYYY-MM-DD"
$ws.Range("G50").Value = 19
$ws.Range("H50").Value = "Major walk report community move green school.`nNorth no peace want term change brother. Education of class job hard example speak minute."
$ws.Range("I50").Value = "Jack Jackson"
$ws.Range("J50").Value = "APROBADO - DECANO"
$ws.Range("K50").Value = "Type 1"
$ws.Range("L50").Value = "Taylor Wilkins (@5464236),Kenneth Chang (@59512309),Melissa Johnson (@47666474)"
